$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("baseline-variables")

# Column A: rename dns1/dns2/dns3 -> tdns1/tdns2/tdns3
$ws.Range("A23").Value = "tdns1"
$ws.Range("A24").Value = "tdns2"
$ws.Range("A25").Value = "tdns3"

# Column B: add full descriptive names
$ws.Range("B23").Value = "Treasury-FFR Spread Level (10-Year Level)"
$ws.Range("B24").Value = "Treasury-FFR Spread Slope (Negative of 10Y-3M Spread)"
$ws.Range("B25").Value = "Treasury-FFR Spread Curvature"

# Column C: category
$ws.Range("C23").Value = "Interest Rates"
$ws.Range("C24").Value = "Interest Rates"
$ws.Range("C25").Value = "Interest Rates"

# Column D: source
$ws.Range("D23").Value = "calc"
$ws.Range("D24").Value = "calc"
$ws.Range("D25").Value = "calc"

# Column F: units
$ws.Range("F23").Value = "%"
$ws.Range("F24").Value = "%"
$ws.Range("F25").Value = "%"

# Column G: freq
$ws.Range("G23").Value = "m"
$ws.Range("G24").Value = "m"
$ws.Range("G25").Value = "m"

# Column I: st
$ws.Range("I23").Value = "d"
$ws.Range("I24").Value = "d"
$ws.Range("I25").Value = "d"

# Column J: d1
$ws.Range("J23").Value = "base"
$ws.Range("J24").Value = "base"
$ws.Range("J25").Value = "base"

# Column K: d2
$ws.Range("K23").Value = "none"
$ws.Range("K24").Value = "none"
$ws.Range("K25").Value = "none"

# Update the active selection to match the saved cursor position
$ws.Activate()
$ws.Range("I4").Select()
